# Apply the header bold formatting across sheets, update the config sheet's
# train_iteration value, and fix the "phase" sheet view (zoom/selection).

$wb = $excel.ActiveWorkbook

$headerSheets = @("phase", "source", "prepare", "train")

foreach ($sheetName in $headerSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $headerRange = $ws.Range("A1:E1")
    $headerRange.Font.Bold = $true
}

# config sheet: bold the "name" column header (A1) and the "value" header (B1)
$configWs = $wb.Worksheets.Item("config")
$configWs.Range("A1").Font.Bold = $true
$configWs.Range("B1").Font.Bold = $true

# config sheet: train_iteration (A6) value 2 -> 20
$configWs.Range("B6").Value = 20

# phase sheet view: reset zoom to 100 (normal) and drop selection override
$phaseWs = $wb.Worksheets.Item("phase")
$phaseWs.Activate()
$excel.ActiveWindow.Zoom = 100
